$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add contact numbers in column B for rows 8, 13, 16
$ws.Range("B8").Value = 9674569343
$ws.Range("B13").Value = 9830304429
$ws.Range("B16").Value = 9836581909

# Update the selected cell to C15
$ws.Range("C15").Select()
